# Update the "Pais" sheet: refresh the COVID country/provincia snapshot
# (row reorders where two countries swap ranking + the updated case counts,
# plus the "last updated" timestamp banner in A1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Banner timestamp ---
$ws.Range("A1").Value = 'Datos actualizados a 22 de Octubre de 2020 a las 19:24'

# --- Updated case counts (countries whose totals changed but rank did not) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 8607785
$ws.Range("C4").Value = 22966
$ws.Range("D4").Value = 5613849
$ws.Range("E4").Value = 2766240
$ws.Range("G4").Value = 287
$ws.Range("H4").Value = 227696

# Row 5: India
$ws.Range("B5").Value = 7756206
$ws.Range("C5").Value = 51048
$ws.Range("D5").Value = 6941238
$ws.Range("E5").Value = 697691
$ws.Range("G5").Value = 624
$ws.Range("H5").Value = 117277

# Row 8: España
$ws.Range("B8").Value = 1090521
$ws.Range("C8").Value = 20986
$ws.Range("G8").Value = 155
$ws.Range("H8").Value = 34521

# Row 14: Reino Unido
$ws.Range("B14").Value = 810467
$ws.Range("C14").Value = 21242
$ws.Range("G14").Value = 189
$ws.Range("H14").Value = 44347

# Row 24: Turquia
$ws.Range("B24").Value = 355528
$ws.Range("C24").Value = 2102
$ws.Range("D24").Value = 310027
$ws.Range("E24").Value = 35917
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 9584

# Row 28: Israel
$ws.Range("B28").Value = 308166
$ws.Range("C28").Value = 831
$ws.Range("D28").Value = 287977
$ws.Range("E28").Value = 17870

# Row 36: Ecuador
$ws.Range("B36").Value = 156451
$ws.Range("C36").Value = 826
$ws.Range("E36").Value = 9764
$ws.Range("G36").Value = 47
$ws.Range("H36").Value = 12500

# Row 90: Republica de Macedonia
$ws.Range("B90").Value = 25473
$ws.Range("C90").Value = 637
$ws.Range("D90").Value = 18047
$ws.Range("E90").Value = 6552
$ws.Range("G90").Value = 12
$ws.Range("H90").Value = 874

# Row 118: Zimbabue
$ws.Range("B118").Value = 8242
$ws.Range("C118").Value = 27
$ws.Range("D118").Value = 7742
$ws.Range("E118").Value = 264

# Row 126: Suazilandia
$ws.Range("B126").Value = 5814
$ws.Range("C126").Value = 9
$ws.Range("D126").Value = 5468
$ws.Range("E126").Value = 230

# Row 154: Sudan del Sur
$ws.Range("B154").Value = 2872
$ws.Range("C154").Value = 2
$ws.Range("E154").Value = 1527

# Row 159: Sierra Leona
$ws.Range("B159").Value = 2340
$ws.Range("C159").Value = 3
$ws.Range("D159").Value = 1777
$ws.Range("E159").Value = 490

# Row 165: Republica del Chad
$ws.Range("B165").Value = 1410
$ws.Range("C165").Value = 6
$ws.Range("D165").Value = 1223
$ws.Range("E165").Value = 91

# --- Ranking swaps: two adjacent countries exchanged places, each row keeps
#     its row number but receives the other country's label + fresh data ---
# Rows 31/32: Chequia <-> Polonia
$ws.Range("A31").Value = 'Chequia'
$ws.Range("B31").Value = 216425
$ws.Range("C31").Value = 7510
$ws.Range("D31").Value = 83649
$ws.Range("E31").Value = 130948
$ws.Range("G31").Value = 89
$ws.Range("H31").Value = 1828
$ws.Range("A32").Value = 'Polonia'
$ws.Range("B32").Value = 214686
$ws.Range("C32").Value = 12107
$ws.Range("D32").Value = 102204
$ws.Range("E32").Value = 108463
$ws.Range("G32").Value = 168
$ws.Range("H32").Value = 4019

# Rows 107/108: Maldivas <-> Mozambique
$ws.Range("A107").Value = 'Maldivas'
$ws.Range("B107").Value = 11358
$ws.Range("C107").Value = 42
$ws.Range("D107").Value = 10383
$ws.Range("E107").Value = 938
$ws.Range("H107").Value = 37
$ws.Range("A108").Value = 'Mozambique'
$ws.Range("B108").Value = 11331
$ws.Range("D108").Value = 9165
$ws.Range("E108").Value = 2087
$ws.Range("H108").Value = 79

# Rows 123/124: Sri Lanka <-> Bahamas
$ws.Range("A123").Value = 'Sri Lanka'
$ws.Range("B123").Value = 6287
$ws.Range("C123").Value = 309
$ws.Range("D123").Value = 3561
$ws.Range("E123").Value = 2712
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 14
$ws.Range("A124").Value = 'Bahamas'
$ws.Range("B124").Value = 6051
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 3633
$ws.Range("E124").Value = 2291
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 127

# Rows 202/203: Santa Lucia <-> Puerto Rico
$ws.Range("A202").Value = 'Santa Lucia'
$ws.Range("B202").Value = 42
$ws.Range("C202").Value = 4
$ws.Range("D202").Value = 27
$ws.Range("E202").Value = 15
$ws.Range("H202").Value = 0
$ws.Range("A203").Value = 'Puerto Rico'
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 1
$ws.Range("E203").Value = 36
